$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the name/balance/transactions columns of the template data row to
# a text format *before* duplicating the row, so the copy below reuses the
# same (now-text) style instead of minting a second, redundant one.
$ws.Range("B4:G4,N4").NumberFormat = "@"
$ws.Range("H4:K4").NumberFormat = "@"

# Insert a new row at position 5. This pushes the old row 5 (the K:N
# subtotal strip) down to row 6 and the old row 6 (footer line) down to
# row 7, automatically re-anchoring their existing merged ranges.
$ws.Rows(5).Insert()

# Duplicate row 4's formatting (fonts/fills/borders/number formats) onto
# the newly inserted row 5 so the second product line looks the same as
# the first.
$ws.Range("A4:N4").Copy()
$ws.Range("A5:N5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 4: first product line ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "ELICA-M CREAM 30 GRAM"
$ws.Range("H4").Value = "0:0"
$ws.Range("L4").Value = 52
$ws.Range("N4").Value = "1:0"

# --- Row 5: second product line ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "VIDROP 2800 I.U./ML ORAL DROPS 15 ML"
$ws.Range("H5").Value = "6:0"
$ws.Range("L5").Value = 26
$ws.Range("N5").Value = "1:0"

# Re-create the merges for row 5 to mirror row 4's layout (PasteSpecial
# formats-only doesn't bring merge state along).
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()

# Row heights: new row 5 matches the data-row height; row 6 (old row 5)
# keeps that same height; row 7 (old row 6, footer) keeps its height.
$ws.Rows(5).RowHeight = 25.5
$ws.Rows(6).RowHeight = 25.5
$ws.Rows(7).RowHeight = 16.5

# --- Row 6: totals strip now carries the sum of the sale-price column ---
$ws.Range("K6").Value = 78
